# Edit script: week 4 day 3
# Adds Pot_weight_g (B), updates Perc_brown (C), adds Diameter_mm (D) data
# to PIPO (sheet1) rows 22-61 and PSME (sheet2) rows 2-61.

$wb = $excel.ActiveWorkbook
$wsPIPO = $wb.Worksheets.Item("PIPO")
$wsPSME = $wb.Worksheets.Item("PSME")

$pipoData = @(
    @(22, 355.2, 0.1, 5.604),
    @(23, 357.8, 0.1, 4.624),
    @(24, 329.1, 0.1, 4.811),
    @(25, 377.7, 0.25, 4.764),
    @(26, 314.0, 0.25, 4.303),
    @(27, 302.4, 0.1, 5.041),
    @(28, 274.8, 0.1, 6.271),
    @(29, 249.0, 0.1, 5.739),
    @(30, 284.5, 0.1, 6.42),
    @(31, 360.1, 0.1, 5.033),
    @(32, 463.8, 0.1, 4.883),
    @(33, 286.4, 0.1, 5.094),
    @(34, 262.7, 0.1, 5.005),
    @(35, 261.7, 0.25, 4.785),
    @(36, 477.2, 0.1, 4.341),
    @(37, 311.5, 0.1, 5.79),
    @(38, 316.1, 0.1, 5.1),
    @(39, 313.9, 0.1, 4.793),
    @(40, 336.0, 0.1, 4.132),
    @(41, 371.5, 0.1, 5.091),
    @(42, 286.8, 0.1, 4.99),
    @(43, 319.4, 0.1, 5.615),
    @(44, 267.0, 0.1, 5.052),
    @(45, 268.5, 0.1, 4.929),
    @(46, 292.8, 0.1, 5.026),
    @(47, 282.5, 0.25, 5.181),
    @(48, 330.2, 0.1, 5.16),
    @(49, 359.6, 0.1, 4.991),
    @(50, 291.0, 0.1, 5.158),
    @(51, 256.9, 0.1, 5.646),
    @(52, 397.5, 0.1, 5.866),
    @(53, 327.0, 0.25, 5.223),
    @(54, 394.1, 0.25, 5.3),
    @(55, 276.3, 0.1, 4.943),
    @(56, 254.7, 0.1, 6.707),
    @(57, 464.5, 0.1, 4.622),
    @(58, 279.3, 0.1, 6.244),
    @(59, 282.4, 0.1, 5.437),
    @(60, 357.3, 0.1, 5.406),
    @(61, 302.2, 0.1, 5.688),
)

$psmeData = @(
    @(2, 311.6, 0.1, 3.366),
    @(3, 349.6, 0.1, 3.701),
    @(4, 509.3, 0.1, 4.016),
    @(5, 336.5, 0.1, 4.089),
    @(6, 409.5, 0.1, 2.678),
    @(7, 380.0, 0.1, 3.061),
    @(8, 299.7, 0.1, 4.658),
    @(9, 307.1, 0.1, 4.477),
    @(10, 323.5, 0.1, 4.082),
    @(11, 251.5, 0.1, 4.076),
    @(12, 288.4, 0.1, 3.799),
    @(13, 485.8, 0.1, 3.99),
    @(14, 401.3, 0.1, 3.251),
    @(15, 359.3, 0.1, 3.761),
    @(16, 267.1, 0.1, 4.376),
    @(17, 319.7, 0.1, 4.143),
    @(18, 460.2, 0.1, 4.142),
    @(19, 390.7, 0.1, 2.914),
    @(20, 369.6, 0.1, 3.25),
    @(21, 295.8, 0.1, 4.362),
    @(22, 296.8, 0.1, 4.112),
    @(23, 539.3, 0.1, 3.344),
    @(24, 361.5, 0.1, 3.725),
    @(25, 377.4, 0.1, 3.811),
    @(26, 315.3, 0.1, 3.287),
    @(27, 247.6, 0.1, 4.601),
    @(28, 515.5, 0.1, 5.029),
    @(29, 272.9, 0.1, 4.297),
    @(30, 340.2, 0.1, 3.504),
    @(31, 279.5, 0.1, 3.635),
    @(32, 345.7, 0.1, 4.942),
    @(33, 314.0, 0.1, 4.439),
    @(34, 334.1, 0.1, 3.584),
    @(35, 359.7, 0.1, 3.361),
    @(36, 316.6, 0.1, 3.11),
    @(37, 282.8, 0.1, 4.183),
    @(38, 402.0, 0.1, 4.831),
    @(39, 289.6, 0.1, 4.468),
    @(40, 283.8, 0.1, 4.361),
    @(41, 254.3, 0.1, 4.849),
    @(42, 255.1, 0.1, 4.597),
    @(43, 281.5, 0.1, 4.615),
    @(44, 266.6, 0.1, 6.236),
    @(45, 248.8, 0.1, 4.722),
    @(46, 422.0, 0.1, 3.919),
    @(47, 334.3, 0.1, 3.698),
    @(48, 236.2, 0.1, 4.994),
    @(49, 421.6, 0.1, 4.33),
    @(50, 270.9, 0.1, 4.313),
    @(51, 244.0, 0.1, 4.7),
    @(52, 305.9, 0.1, 3.891),
    @(53, 349.9, 0.1, 4.0),
    @(54, 304.8, 0.1, 3.848),
    @(55, 255.7, 0.1, 4.493),
    @(56, 250.3, 0.1, 4.966),
    @(57, 269.8, 0.1, 4.27),
    @(58, 354.5, 0.1, 4.246),
    @(59, 278.4, 0.1, 4.584),
    @(60, 259.8, 0.1, 6.283),
    @(61, 321.6, 0.1, 4.383),
)

foreach ($row in $pipoData) {
    $r = $row[0]
    $wsPIPO.Cells.Item($r, 2).Value = $row[1]
    $wsPIPO.Cells.Item($r, 3).Value = $row[2]
    $wsPIPO.Cells.Item($r, 4).Value = $row[3]
}

foreach ($row in $psmeData) {
    $r = $row[0]
    $wsPSME.Cells.Item($r, 2).Value = $row[1]
    $wsPSME.Cells.Item($r, 3).Value = $row[2]
    $wsPSME.Cells.Item($r, 4).Value = $row[3]
}

# Update PSME data validation range from C2:C31 to C2:C61
$psmeValidation = $wsPSME.Range("C2:C61")
$psmeValidation.Validation.Delete()
$psmeValidation.Validation.Add(3, 1, 1, '"10%, 25%, 50%, 75%, 90%"')

# Update sheet view / selection state for PIPO
$wsPIPO.Activate()
$excel.ActiveWindow.ScrollRow = 46
$wsPIPO.Range("D62").Select()

# Update sheet view / selection / freeze panes state for PSME
$wsPSME.Activate()
$excel.ActiveWindow.Zoom = 174
$wsPSME.Range("A43").Select()
$excel.ActiveWindow.FreezePanes = $false
$wsPSME.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsPSME.Range("C62").Select()

